# Insert a new weekly price record as row 569 in the "Acelga" (Hortaliza)
# price sheet, pushing the existing rows 569-633 down to 570-634.
#
# This mirrors the author's edit: a new observation dated 45212 (2023-10-13)
# was added at the top of the data block, shifting every following row down
# by one, with the former last row (633) becoming the new last row (634).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 569:633 down to 570:634, leaving row 569 blank and ready
# to receive the new record.
$ws.Rows.Item(569).Insert()

# Populate the newly inserted row with the new weekly record.
$ws.Cells.Item(569, 1).Value2 = 3
$ws.Cells.Item(569, 2).Value2 = "Femacal de La Calera"
$ws.Cells.Item(569, 3).Value2 = "Coquimbo"
$ws.Cells.Item(569, 4).Value2 = 45212
$ws.Cells.Item(569, 5).Value2 = 5
$ws.Cells.Item(569, 6).Value2 = 100112009
$ws.Cells.Item(569, 7).Value2 = "Acelga"
$ws.Cells.Item(569, 8).Value2 = "Sin especificar"
$ws.Cells.Item(569, 9).Value2 = "Primera"
$ws.Cells.Item(569, 10).Value2 = 230
$ws.Cells.Item(569, 11).Value2 = 3000
$ws.Cells.Item(569, 12).Value2 = 3500
$ws.Cells.Item(569, 13).Value2 = 3261
$ws.Cells.Item(569, 14).Value2 = "$/docena de atados (6 kilos)"
$ws.Cells.Item(569, 15).Value2 = "Provincia de Quillota"
$ws.Cells.Item(569, 16).Value2 = 544
$ws.Cells.Item(569, 17).Value2 = 6
$ws.Cells.Item(569, 18).Value2 = "Hortaliza"
